$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 4).Value = '''27.002.24'
$ws.Cells.Item(2, 5).Value = '  +5.29%  '

# Row 3
$ws.Cells.Item(3, 4).Value = '''1.883.10'
$ws.Cells.Item(3, 5).Value = '  +4.27%  '

# Row 4
$ws.Cells.Item(4, 4).Value = '''1.001'
$ws.Cells.Item(4, 5).Value = '  +0.07%  '

# Row 5
$ws.Cells.Item(5, 4).Value = '''283.11'
$ws.Cells.Item(5, 5).Value = '  +2.40%  '

# Row 6
$ws.Cells.Item(6, 4).Value = '''1.000'
$ws.Cells.Item(6, 5).Value = '  +0.03%  '

# Row 7
$ws.Cells.Item(7, 4).Value = '''0.5273'
$ws.Cells.Item(7, 5).Value = '  +4.11%  '

# Row 8
$ws.Cells.Item(8, 4).Value = '''0.3544'
$ws.Cells.Item(8, 5).Value = '  +0.86%  '

# Row 9
$ws.Cells.Item(9, 4).Value = '''45.30'
$ws.Cells.Item(9, 5).Value = '  +3.94%  '

# Row 10
$ws.Cells.Item(10, 4).Value = '''0.07097'
$ws.Cells.Item(10, 5).Value = '  +6.18%  '

# Row 11
$ws.Cells.Item(11, 4).Value = '''20.43'
$ws.Cells.Item(11, 5).Value = '  +2.00%  '

# Row 12
$ws.Cells.Item(12, 4).Value = '''0.8205'
$ws.Cells.Item(12, 5).Value = '  -1.96%  '

# Row 13
$ws.Cells.Item(13, 4).Value = '''0.07823'
$ws.Cells.Item(13, 5).Value = '  +0.53%  '

# Row 14
$ws.Cells.Item(14, 4).Value = '''1.892.82'
$ws.Cells.Item(14, 5).Value = '  +4.76%  '

# Row 15
$ws.Cells.Item(15, 4).Value = '''5.233'
$ws.Cells.Item(15, 5).Value = '  +3.26%  '

# Row 16
$ws.Cells.Item(16, 4).Value = '''90.71'
$ws.Cells.Item(16, 5).Value = '  +3.77%  '

# Row 17
$ws.Cells.Item(17, 4).Value = '''1.002'
$ws.Cells.Item(17, 5).Value = '  +0.25%  '

# Row 18
$ws.Cells.Item(18, 4).Value = '''14.57'
$ws.Cells.Item(18, 5).Value = '  +4.67%  '

# Row 19
$ws.Cells.Item(19, 4).Value = '''0.000008172'
$ws.Cells.Item(19, 5).Value = '  +2.65%  '

# Row 20
$ws.Cells.Item(20, 4).Value = '''0.9994'
$ws.Cells.Item(20, 5).Value = '  -0.04%  '

# Row 21
$ws.Cells.Item(21, 4).Value = '''27.056.34'
$ws.Cells.Item(21, 5).Value = '  +5.26%  '

# Row 22
$ws.Cells.Item(22, 2).Value = 'Uniswap'
$ws.Cells.Item(22, 3).Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Cells.Item(22, 4).Value = '''4.799'
$ws.Cells.Item(22, 5).Value = '  +1.80%  '

# Row 23
$ws.Cells.Item(23, 2).Value = 'Cosmos'
$ws.Cells.Item(23, 3).Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Cells.Item(23, 4).Value = '''10.22'
$ws.Cells.Item(23, 5).Value = '  +2.16%  '

# Row 24
$ws.Cells.Item(24, 2).Value = 'Chainlink'
$ws.Cells.Item(24, 3).Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Cells.Item(24, 4).Value = '''6.280'
$ws.Cells.Item(24, 5).Value = '  +4.06%  '

# Row 25
$ws.Cells.Item(25, 2).Value = 'LidoDAOToken'
$ws.Cells.Item(25, 3).Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Cells.Item(25, 4).Value = '''2.416'
$ws.Cells.Item(25, 5).Value = '  +14.42%  '

# Row 26
$ws.Cells.Item(26, 2).Value = 'Monero'
$ws.Cells.Item(26, 3).Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Cells.Item(26, 4).Value = '''147.28'
$ws.Cells.Item(26, 5).Value = '  +3.60%  '

# Row 27
$ws.Cells.Item(27, 2).Value = 'EthereumClassic'
$ws.Cells.Item(27, 3).Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Cells.Item(27, 4).Value = '''17.62'
$ws.Cells.Item(27, 5).Value = '  +4.28%  '

# Row 28
$ws.Cells.Item(28, 2).Value = 'Toncoin'
$ws.Cells.Item(28, 3).Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Cells.Item(28, 4).Value = '''1.673'
$ws.Cells.Item(28, 5).Value = '  +1.04%  '

# Row 29
$ws.Cells.Item(29, 2).Value = 'BitcoinCash'
$ws.Cells.Item(29, 3).Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Cells.Item(29, 4).Value = '''113.50'
$ws.Cells.Item(29, 5).Value = '  +4.58%  '

# Row 30
$ws.Cells.Item(30, 2).Value = 'InternetComputer(DFINITY)'
$ws.Cells.Item(30, 3).Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Cells.Item(30, 4).Value = '''4.419'
$ws.Cells.Item(30, 5).Value = '  +2.41%  '

# Row 31
$ws.Cells.Item(31, 2).Value = 'Filecoin'
$ws.Cells.Item(31, 3).Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Cells.Item(31, 4).Value = '''4.391'
$ws.Cells.Item(31, 5).Value = '  +3.95%  '

# Row 32
$ws.Cells.Item(32, 2).Value = 'Stellar'
$ws.Cells.Item(32, 3).Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Cells.Item(32, 4).Value = '''0.08882'
$ws.Cells.Item(32, 5).Value = '  +0.78%  '

# Row 33
$ws.Cells.Item(33, 2).Value = 'Hedera'
$ws.Cells.Item(33, 3).Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Cells.Item(33, 4).Value = '''0.04923'
$ws.Cells.Item(33, 5).Value = '  +2.81%  '

# Row 34
$ws.Cells.Item(34, 2).Value = 'ARBITRUM'
$ws.Cells.Item(34, 3).Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Cells.Item(34, 4).Value = '''1.180'
$ws.Cells.Item(34, 5).Value = '  +4.87%  '

# Row 35
$ws.Cells.Item(35, 2).Value = 'ImmutableX'
$ws.Cells.Item(35, 3).Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Cells.Item(35, 4).Value = '''0.7484'
$ws.Cells.Item(35, 5).Value = '  +3.25%  '

# Row 36
$ws.Cells.Item(36, 2).Value = 'HuobiToken'
$ws.Cells.Item(36, 3).Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Cells.Item(36, 4).Value = '''2.900'
$ws.Cells.Item(36, 5).Value = '  +1.81%  '

# Row 37
$ws.Cells.Item(37, 2).Value = 'MXToken'
$ws.Cells.Item(37, 3).Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Cells.Item(37, 4).Value = '''3.302'
$ws.Cells.Item(37, 5).Value = '  +8.85%  '

# Row 38
$ws.Cells.Item(38, 2).Value = 'RenderToken'
$ws.Cells.Item(38, 3).Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Cells.Item(38, 4).Value = '''2.404'
$ws.Cells.Item(38, 5).Value = '  +4.29%  '

# Row 39
$ws.Cells.Item(39, 2).Value = 'TheSandbox'
$ws.Cells.Item(39, 3).Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Cells.Item(39, 4).Value = '''0.5328'
$ws.Cells.Item(39, 5).Value = '  +3.27%  '

# Row 40
$ws.Cells.Item(40, 2).Value = 'VeChain'
$ws.Cells.Item(40, 3).Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Cells.Item(40, 4).Value = '''0.01892'
$ws.Cells.Item(40, 5).Value = '  +1.94%  '

# Row 41
$ws.Cells.Item(41, 2).Value = 'TrustWalletToken'
$ws.Cells.Item(41, 3).Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Cells.Item(41, 4).Value = '''0.9871'
$ws.Cells.Item(41, 5).Value = '  +2.55%  '

# Row 42
$ws.Cells.Item(42, 2).Value = 'Quant'
$ws.Cells.Item(42, 3).Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Cells.Item(42, 4).Value = '''117.09'
$ws.Cells.Item(42, 5).Value = '  +2.23%  '

# Row 43
$ws.Cells.Item(43, 2).Value = 'FraxShare'
$ws.Cells.Item(43, 3).Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Cells.Item(43, 4).Value = '''6.313'
$ws.Cells.Item(43, 5).Value = '  +2.05%  '

# Row 44
$ws.Cells.Item(44, 2).Value = 'Aptos'
$ws.Cells.Item(44, 3).Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Cells.Item(44, 4).Value = '''8.216'
$ws.Cells.Item(44, 5).Value = '  +2.17%  '

# Row 45
$ws.Cells.Item(45, 2).Value = 'Decentraland'
$ws.Cells.Item(45, 3).Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Cells.Item(45, 4).Value = '''0.4634'
$ws.Cells.Item(45, 5).Value = '  +1.00%  '

# Row 46
$ws.Cells.Item(46, 2).Value = 'PaxDollar'
$ws.Cells.Item(46, 3).Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Cells.Item(46, 4).Value = '''0.9997'
$ws.Cells.Item(46, 5).Value = '  +0.03%  '

# Row 47
$ws.Cells.Item(47, 2).Value = 'Algorand'
$ws.Cells.Item(47, 3).Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Cells.Item(47, 4).Value = '''0.1373'
$ws.Cells.Item(47, 5).Value = '  -1.04%  '

# Row 48
$ws.Cells.Item(48, 2).Value = 'EnergySwap'
$ws.Cells.Item(48, 3).Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Cells.Item(48, 4).Value = '''9.546'
$ws.Cells.Item(48, 5).Value = '  +3.30%  '

# Row 49
$ws.Cells.Item(49, 2).Value = 'Elrond'
$ws.Cells.Item(49, 3).Value = 'https://coinranking.com/coin/omwkOTglq+elrond-egld'
$ws.Cells.Item(49, 4).Value = '''36.84'
$ws.Cells.Item(49, 5).Value = '  +2.86%  '

# Row 50
$ws.Cells.Item(50, 2).Value = 'NEARProtocol'
$ws.Cells.Item(50, 3).Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Cells.Item(50, 4).Value = '''1.530'
$ws.Cells.Item(50, 5).Value = '  +2.68%  '

# Row 51
$ws.Cells.Item(51, 2).Value = 'Cronos'
$ws.Cells.Item(51, 3).Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Cells.Item(51, 4).Value = '''0.05951'
$ws.Cells.Item(51, 5).Value = '  +2.34%  '
